$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "MD410 Attendance": refresh timestamp in title, and flip two voters'
# "Voter" flags (E15: Yes -> No, E24: No -> Yes)
# ---------------------------------------------------------------------------
$wsMD410 = $wb.Worksheets.Item("MD410 Attendance")
$wsMD410.Range("A1").Value = "MD410 Registrees as of 21/04/2021 12:20"
$wsMD410.Cells.Item(15, 5).Value = "No"
$wsMD410.Cells.Item(24, 5).Value = "Yes"

# ---------------------------------------------------------------------------
# Sheet "410E Attendance": refresh timestamp, flip several "Voter" flags
# ---------------------------------------------------------------------------
$ws410E = $wb.Worksheets.Item("410E Attendance")
$ws410E.Range("A1").Value = "410E Registrees as of 21/04/2021 12:20"
$ws410E.Cells.Item(8, 5).Value = "No"
$ws410E.Cells.Item(9, 5).Value = "No"
$ws410E.Cells.Item(13, 5).Value = "Yes"
$ws410E.Cells.Item(91, 5).Value = "No"
$ws410E.Cells.Item(99, 5).Value = "Yes"
$ws410E.Cells.Item(117, 5).Value = "Yes"

# ---------------------------------------------------------------------------
# Sheet "410W Attendance": refresh timestamp only
# ---------------------------------------------------------------------------
$ws410W = $wb.Worksheets.Item("410W Attendance")
$ws410W.Range("A1").Value = "410W Registrees as of 21/04/2021 12:20"

# ---------------------------------------------------------------------------
# Sheet "410E Voting": refresh timestamp, add two new club rows
# ("East Coast" and "Gonubie"), bump Benoni Lakes's vote count, and refresh
# the trailing summary rows.
# ---------------------------------------------------------------------------
$ws410EV = $wb.Worksheets.Item("410E Voting")
$ws410EV.Range("A1").Value = "410E Voting details as of 21/04/2021 12:20"

# Insert "East Coast" just above "East London Beacon Bay" (row 7)
$ws410EV.Rows.Item(7).Insert()
$ws410EV.Range("A7:B7").Borders.LineStyle = 1
$ws410EV.Rows.Item(7).RowHeight = 25
$ws410EV.Cells.Item(7, 1).Value = "East Coast"
$ws410EV.Cells.Item(7, 2).Value = 1

# Insert "Gonubie" just above "Helderkruin" (now at row 11 after the above insert)
$ws410EV.Rows.Item(11).Insert()
$ws410EV.Range("A11:B11").Borders.LineStyle = 1
$ws410EV.Rows.Item(11).RowHeight = 25
$ws410EV.Cells.Item(11, 1).Value = "Gonubie"
$ws410EV.Cells.Item(11, 2).Value = 1

# Benoni Lakes now has 2 voters instead of 1
$ws410EV.Cells.Item(4, 2).Value = 2

# Trailing summary rows shifted down from 32/33 to 34/35, with updated club count
$ws410EV.Cells.Item(34, 1).Value = "Number of clubs: 31"
$ws410EV.Cells.Item(35, 1).Value = "Number of voters: 52"

# ---------------------------------------------------------------------------
# Sheet "410W Voting": refresh timestamp only
# ---------------------------------------------------------------------------
$ws410WV = $wb.Worksheets.Item("410W Voting")
$ws410WV.Range("A1").Value = "410W Voting details as of 21/04/2021 12:20"
